# Normalize cell text: replace embedded line breaks with a single space.
# (Mirrors an upload/re-save that flattened multi-line cell text onto one line.)

$wb = $excel.ActiveWorkbook
$lf = [char]10

foreach ($ws in $wb.Worksheets) {
    $used = $ws.UsedRange
    $rowCount = $used.Rows.Count
    $colCount = $used.Columns.Count
    $startRow = $used.Row
    $startCol = $used.Column

    for ($r = 0; $r -lt $rowCount; $r++) {
        for ($c = 0; $c -lt $colCount; $c++) {
            $cell = $ws.Cells.Item($startRow + $r, $startCol + $c)
            $v = $cell.Value2

            if (($v -ne $null) -and ($v -is [string]) -and $v.Contains($lf)) {
                $new = $v.Replace($lf, " ")
                $cell.Value = $new
            }
        }
    }
}
